# feat: add 2022-Q3 data
#
# Layout before:
#   Sheet1 "总计"     (sheetId=1, rId1)
#   Sheet2 "2021-Q4"  (sheetId=2, rId2)  <- fund holdings for 2021-Q4
#
# Layout after:
#   Sheet1 "总计"     (sheetId=1, rId1)  <- gains a new summary row for 2022-Q3
#   Sheet2 "2022-Q3"  (sheetId=2, rId2)  <- NEW fund holdings for 2022-Q3
#   Sheet3 "2021-Q4"  (sheetId=3, rId3)  <- old fund holdings, unchanged, just shifted over

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item(1)
$oldQ4Sheet = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# 1) Turn the existing "2021-Q4" sheet into the new "2022-Q3" sheet in place
#    (this keeps sheetId=2 / rId2, matching the target numbering), then add a
#    fresh sheet right after it to hold the old 2021-Q4 data (picks up
#    sheetId=3 / rId3).
# ---------------------------------------------------------------------------
$oldQ4Sheet.Name = "2022-Q3"
$newQ4Sheet = $wb.Worksheets.Add($null, $oldQ4Sheet)
$newQ4Sheet.Name = "2021-Q4"

# Copy the untouched 2021-Q4 fund-holdings table (all formatting + values)
# from the renamed sheet over to the brand-new "2021-Q4" sheet, then clear it
# off the (soon to be) "2022-Q3" sheet so we can put fresh data there.
$oldQ4Sheet.Range("A1:H3").Copy()
$newQ4Sheet.Range("A1:H3").PasteSpecial(-4104)
$oldQ4Sheet.Range("A1:H3").ClearContents()

# ---------------------------------------------------------------------------
# 2) Populate "2022-Q3" with the new fund holdings data. Header row + layout
#    mirrors the other quarter sheets.
# ---------------------------------------------------------------------------
$q3 = $oldQ4Sheet

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "'014155"
$q3.Range("C2").Value = "国泰君安中证500指数增强A"
$q3.Range("D2").Value = "'6.64"
$q3.Range("E2").Value = "'92.15"
$q3.Range("F2").Value = "'1.20"
$q3.Range("G2").Value = "'0.0797"
$q3.Range("H2").Value = 7

$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "'014156"
$q3.Range("C3").Value = "国泰君安中证500指数增强C"
$q3.Range("D3").Value = "'4.02"
$q3.Range("E3").Value = "'92.15"
$q3.Range("F3").Value = "'1.20"
$q3.Range("G3").Value = "'0.0482"
$q3.Range("H3").Value = 7

# Re-apply the header / first-column formatting (bold, centered, bordered)
# that the rest of the workbook's tables use, by copying it across from the
# "2021-Q4" sheet we just populated.
$newQ4Sheet.Range("B1:H1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)
$newQ4Sheet.Range("A2:A3").Copy()
$q3.Range("A2:A3").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3) Update the "总计" summary sheet: existing row 2 becomes the 2022-Q3
#    entry, and a new row 3 carries the 2021-Q4 entry that used to live in
#    row 2.
# ---------------------------------------------------------------------------
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q4"
$totalSheet.Range("C3").Value = 2
$totalSheet.Range("D3").Value = 0.12

$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.13

$totalSheet.Range("A2").Copy()
$totalSheet.Range("A3").PasteSpecial(-4122)
